# carrieres.xlsx — "Permet le lancement avec N carrieres"
# * Adds a new "debut" sheet holding PlafondSS data (COR1..COR4) for
#   generations 1930-1990 (used for COR2 init).
# * Adds a "PlafondSS" row to the "meta" sheet.
# * Re-points the active sheet/selection to reflect the edit session.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "meta" sheet: document the new "PlafondSS" series (row 8)
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("meta")
$meta.Cells.Item(8, 1).Value = "PlafondSS"
$meta.Cells.Item(8, 2).Value = "PlafondSS"
$meta.Cells.Item(8, 3).Value = "PlafondSS"

# ---------------------------------------------------------------------
# 2. New "debut" sheet, inserted right after "serie"
# ---------------------------------------------------------------------
$serie = $wb.Worksheets.Item("serie")
$debut = $wb.Worksheets.Add($null, $serie)
$debut.Name = "debut"

$debut.Cells.Item(1, 1).Value = "generation"
$debut.Cells.Item(1, 2).Value = "COR1"
$debut.Cells.Item(1, 3).Value = "COR2"
$debut.Cells.Item(1, 4).Value = "COR3"
$debut.Cells.Item(1, 5).Value = "COR4"

$debutData = @(
    @(1930, 19.525, 17.13, 17.13, 19.67),
    @(1931, 19.525, 17.13, 17.13, 19.67),
    @(1932, 19.525, 17.13, 17.13, 19.67),
    @(1933, 19.525, 17.13, 17.13, 19.67),
    @(1934, 19.525, 17.13, 17.13, 19.67),
    @(1935, 19.525, 17.13, 17.13, 19.67),
    @(1936, 19.525, 17.13, 17.13, 19.67),
    @(1937, 19.525, 17.13, 17.13, 19.67),
    @(1938, 19.525, 17.13, 17.13, 19.67),
    @(1939, 19.525, 17.13, 17.13, 19.67),
    @(1940, 19.525, 17.13, 17.13, 19.67),
    @(1941, 19.525, 17.13, 17.13, 19.67),
    @(1942, 19.525, 17.13, 17.13, 19.67),
    @(1943, 19.525, 17.13, 17.13, 19.67),
    @(1944, 19.525, 17.13, 17.13, 19.67),
    @(1945, 19.525, 17.13, 17.13, 19.67),
    @(1946, 19.525, 17.13, 17.13, 19.67),
    @(1947, 19.341875, 17.179375, 17.179375, 19.5725),
    @(1948, 19.15875, 17.22875, 17.22875, 19.475),
    @(1949, 18.975625, 17.278125, 17.278125, 19.3775),
    @(1950, 18.7925, 17.3275, 17.3275, 19.28),
    @(1951, 18.961875, 17.4425, 17.4425, 19.38875),
    @(1952, 19.13125, 17.5575, 17.5575, 19.4975),
    @(1953, 19.300625, 17.6725, 17.6725, 19.60625),
    @(1954, 19.47, 17.7875, 17.7875, 19.715),
    @(1955, 19.57625, 17.898125, 17.898125, 19.81875),
    @(1956, 19.6825, 18.00875, 18.00875, 19.9225),
    @(1957, 19.78875, 18.119375, 18.119375, 20.02625),
    @(1958, 19.895, 18.23, 18.23, 20.13),
    @(1959, 19.99, 18.328125, 18.328125, 20.249375),
    @(1960, 20.085, 18.42625, 18.42625, 20.36875),
    @(1961, 20.18, 18.524375, 18.524375, 20.488125),
    @(1962, 20.275, 18.6225, 18.6225, 20.6075),
    @(1963, 20.40125, 18.82625, 18.82625, 20.746875),
    @(1964, 20.5275, 19.03, 19.03, 20.88625),
    @(1965, 20.65375, 19.23375, 19.23375, 21.025625),
    @(1966, 20.78, 19.4375, 19.4375, 21.165),
    @(1967, 20.90625, 19.533125, 19.533125, 21.259375),
    @(1968, 21.0325, 19.62875, 19.62875, 21.35375),
    @(1969, 21.15875, 19.724375, 19.724375, 21.448125),
    @(1970, 21.285, 19.82, 19.82, 21.5425),
    @(1971, 21.374375, 19.990625, 19.990625, 21.656875),
    @(1972, 21.46375, 20.16125, 20.16125, 21.77125),
    @(1973, 21.553125, 20.331875, 20.331875, 21.885625),
    @(1974, 21.6425, 20.5025, 20.5025, 22),
    @(1975, 21.645, 20.538125, 20.538125, 21.95375),
    @(1976, 21.6475, 20.57375, 20.57375, 21.9075),
    @(1977, 21.65, 20.609375, 20.609375, 21.86125),
    @(1978, 21.6525, 20.645, 20.645, 21.815),
    @(1979, 21.6525, 20.645, 20.645, 21.815),
    @(1980, 21.6525, 20.645, 20.645, 21.815),
    @(1981, 21.6525, 20.645, 20.645, 21.815),
    @(1982, 21.6525, 20.645, 20.645, 21.815),
    @(1983, 21.6525, 20.645, 20.645, 21.815),
    @(1984, 21.6525, 20.645, 20.645, 21.815),
    @(1985, 21.6525, 20.645, 20.645, 21.815),
    @(1986, 21.6525, 20.645, 20.645, 21.815),
    @(1987, 21.6525, 20.645, 20.645, 21.815),
    @(1988, 21.6525, 20.645, 20.645, 21.815),
    @(1989, 21.6525, 20.645, 20.645, 21.815),
    @(1990, 21.6525, 20.645, 20.645, 21.815)
)

for ($i = 0; $i -lt $debutData.Length; $i++) {
    $row = $debutData[$i]
    $r = $i + 2
    for ($col = 1; $col -le 5; $col++) {
        $debut.Cells.Item($r, $col).Value = $row[$col - 1]
    }
}

[void]$debut.Range("C12").Select()
$excel.ActiveWindow.Zoom = 180

# ---------------------------------------------------------------------
# 3. Selection / active-sheet bookkeeping
#    "serie" cursor moves to B1 and loses focus; "meta" becomes the
#    active tab with its cursor on the freshly-added C8 cell.
# ---------------------------------------------------------------------
[void]$serie.Range("B1").Select()

[void]$meta.Activate()
[void]$meta.Range("C8").Select()
